# Update NATMI LR-pair TPM-derived statistics (Epo-Ephb4) with recomputed
# values after the underlying TPM script was updated ("update scripts
# wuth new tpm"). Only the numeric result columns (E..T, excluding the
# constant K/L columns) change; identifier columns A-D, K, L stay as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05285566666666667
$ws.Range("H2").Value = 0.158567
$ws.Range("I2").Value = 0.5476571215423245
$ws.Range("J2").Value = 0.5476571215423245
$ws.Range("M2").Value = 25.11140833333333
$ws.Range("N2").Value = 75.334225
$ws.Range("O2").Value = 0.7431105026796001
$ws.Range("P2").Value = 0.7431105026796001
$ws.Range("Q2").Value = 1.327280228397222
$ws.Range("R2").Value = 11.945522055575
$ws.Range("S2").Value = 0.4069697588853796
$ws.Range("T2").Value = 0.4069697588853796

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05285566666666667
$ws.Range("H3").Value = 0.158567
$ws.Range("I3").Value = 0.5476571215423245
$ws.Range("J3").Value = 0.5476571215423245
$ws.Range("O3").Value = 0.1596166092346045
$ws.Range("P3").Value = 0.1596166092346045
$ws.Range("Q3").Value = 0.2850934938975556
$ws.Range("R3").Value = 2.565841445078
$ws.Range("S3").Value = 0.08741517276376949
$ws.Range("T3").Value = 0.0874151727637695

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05285566666666667
$ws.Range("H4").Value = 0.158567
$ws.Range("I4").Value = 0.5476571215423245
$ws.Range("J4").Value = 0.5476571215423245
$ws.Range("N4").Value = 9.861222
$ws.Range("O4").Value = 0.09727288808579543
$ws.Range("P4").Value = 0.09727288808579541
$ws.Range("Q4").Value = 0.1737404876526667
$ws.Range("R4").Value = 1.563664388874
$ws.Range("S4").Value = 0.05327218989317539
$ws.Range("T4").Value = 0.05327218989317539

$ws.Range("I5").Value = 0.3108894545429426
$ws.Range("J5").Value = 0.3108894545429427
$ws.Range("M5").Value = 25.11140833333333
$ws.Range("N5").Value = 75.334225
$ws.Range("O5").Value = 0.7431105026796001
$ws.Range("P5").Value = 0.7431105026796001
$ws.Range("Q5").Value = 0.7534594365722221
$ws.Range("R5").Value = 6.78113492915
$ws.Range("S5").Value = 0.2310252188431928
$ws.Range("T5").Value = 0.2310252188431928

$ws.Range("I6").Value = 0.3108894545429426
$ws.Range("J6").Value = 0.3108894545429427
$ws.Range("O6").Value = 0.1596166092346045
$ws.Range("P6").Value = 0.1596166092346045
$ws.Range("S6").Value = 0.04962312058094021
$ws.Range("T6").Value = 0.04962312058094022

$ws.Range("I7").Value = 0.3108894545429426
$ws.Range("J7").Value = 0.3108894545429427
$ws.Range("N7").Value = 9.861222
$ws.Range("O7").Value = 0.09727288808579543
$ws.Range("P7").Value = 0.09727288808579541
$ws.Range("Q7").Value = 0.09862755967866667
$ws.Range("R7").Value = 0.8876480371079999
$ws.Range("S7").Value = 0.03024111511880965
$ws.Range("T7").Value = 0.03024111511880965

$ws.Range("I8").Value = 0.1414534239147328
$ws.Range("J8").Value = 0.1414534239147328
$ws.Range("M8").Value = 25.11140833333333
$ws.Range("N8").Value = 75.334225
$ws.Range("O8").Value = 0.7431105026796001
$ws.Range("P8").Value = 0.7431105026796001
$ws.Range("Q8").Value = 0.3428209465666667
$ws.Range("R8").Value = 3.0853885191
$ws.Range("S8").Value = 0.1051155249510277
$ws.Range("T8").Value = 0.1051155249510277

$ws.Range("I9").Value = 0.1414534239147328
$ws.Range("J9").Value = 0.1414534239147328
$ws.Range("O9").Value = 0.1596166092346045
$ws.Range("P9").Value = 0.1596166092346045
$ws.Range("S9").Value = 0.02257831588989476
$ws.Range("T9").Value = 0.02257831588989477

$ws.Range("I10").Value = 0.1414534239147328
$ws.Range("J10").Value = 0.1414534239147328
$ws.Range("N10").Value = 9.861222
$ws.Range("O10").Value = 0.09727288808579543
$ws.Range("P10").Value = 0.09727288808579541
$ws.Range("Q10").Value = 0.044875134248
$ws.Range("S10").Value = 0.01375958307381038
$ws.Range("T10").Value = 0.01375958307381038
